$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.474.76"
$ws.Range("E2").Value = "  +2.47%  "
# Row 3
$ws.Range("D3").Value = "3.435.78"
$ws.Range("E3").Value = "  +3.88%  "
# Row 4
$ws.Range("E4").Value = "  -0.32%  "
# Row 5
$ws.Range("E5").Value = "  +3.88%  "
# Row 6
$ws.Range("D6").Value = "'130.44"
$ws.Range("E6").Value = "  +6.07%  "
# Row 7
$ws.Range("D7").Value = "'0.602"
$ws.Range("E7").Value = "  +3.95%  "
# Row 8
$ws.Range("E8").Value = "  -0.04%  "
# Row 9
$ws.Range("D9").Value = "'0.698"
$ws.Range("E9").Value = "  +8.89%  "
# Row 10
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "  +24.59%  "
# Row 11
$ws.Range("D11").Value = "'42.33"
$ws.Range("E11").Value = "  +6.47%  "
# Row 12
$ws.Range("E12").Value = "  +0.76%  "
# Row 13
$ws.Range("D13").Value = "'8.57"
$ws.Range("E13").Value = "  +6.50%  "
# Row 14
$ws.Range("D14").Value = "'19.98"
$ws.Range("E14").Value = "  +6.10%  "
# Row 15
$ws.Range("D15").Value = "3.431.01"
$ws.Range("E15").Value = "  +2.99%  "
# Row 16
$ws.Range("D16").Value = "62.501.62"
$ws.Range("E16").Value = "  +2.41%  "
# Row 17
$ws.Range("D17").Value = "'11.61"
$ws.Range("E17").Value = "  +6.80%  "
# Row 18
$ws.Range("D18").Value = "'0.0000169"
$ws.Range("E18").Value = "  +38.37%  "
# Row 19
$ws.Range("E19").Value = "  +4.26%  "
# Row 20
$ws.Range("D20").Value = "'3.20"
$ws.Range("E20").Value = "  +2.31%  "
# Row 21
$ws.Range("D21").Value = "'84.93"
$ws.Range("E21").Value = "  +8.62%  "
# Row 22
$ws.Range("D22").Value = "'314.94"
$ws.Range("E22").Value = "  +7.94%  "
# Row 23
$ws.Range("E23").Value = "  +3.97%  "
# Row 24
$ws.Range("E24").Value = "  +2.82%  "
# Row 25
$ws.Range("E25").Value = "  +1.86%  "
# Row 26
$ws.Range("D26").Value = "'30.02"
$ws.Range("E26").Value = "  +6.35%  "
# Row 27
$ws.Range("D27").Value = "'8.20"
$ws.Range("E27").Value = "  +3.72%  "
# Row 28
$ws.Range("D28").Value = "'7.81"
$ws.Range("E28").Value = "  +7.32%  "
# Row 29
$ws.Range("D29").Value = "'2.72"
$ws.Range("E29").Value = "  +10.43%  "
# Row 30
$ws.Range("E30").Value = "  +3.04%  "
# Row 31
$ws.Range("D31").Value = "'44.37"
$ws.Range("E31").Value = "  +10.94%  "
# Row 32
$ws.Range("E32").Value = "  +4.34%  "
# Row 33
$ws.Range("D33").Value = "'11.49"
$ws.Range("E33").Value = "  +4.27%  "
# Row 34
$ws.Range("E34").Value = "  +0.13%  "
# Row 35
$ws.Range("D35").Value = "'0.0488"
$ws.Range("E35").Value = "  +5.26%  "
# Row 36
$ws.Range("D36").Value = "'51.45"
$ws.Range("E36").Value = "  -0.52%  "
# Row 37
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.30%  "
# Row 38
$ws.Range("E38").Value = "  +5.46%  "
# Row 39
$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +1.40%  "
# Row 40
$ws.Range("D40").Value = "'0.321"
$ws.Range("E40").Value = "  +18.28%  "
# Row 41
$ws.Range("D41").Value = "'143.97"
$ws.Range("E41").Value = "  +6.50%  "
# Row 42
$ws.Range("E42").Value = "  +5.69%  "
# Row 43
$ws.Range("E43").Value = "  +4.24%  "
# Row 44
$ws.Range("D44").Value = "'17.00"
$ws.Range("E44").Value = "  +5.24%  "
# Row 45
$ws.Range("E45").Value = "  +5.33%  "
# Row 46
$ws.Range("E46").Value = "  +0.64%  "
# Row 47
$ws.Range("D47").Value = "'21.40"
$ws.Range("E47").Value = "  +3.46%  "
# Row 48
$ws.Range("D48").Value = "2.109.05"
$ws.Range("E48").Value = "  +2.01%  "
# Row 49
$ws.Range("E49").Value = "  +12.88%  "
# Row 50
$ws.Range("E50").Value = "  +0.66%  "
# Row 51
$ws.Range("D51").Value = "'1.09"
$ws.Range("E51").Value = "  +34.40%  "
